$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @("PWP201", 1, 0.00000001, 1, 10, 1),
    @("KC200GT2", 1, 0.00000001, 1, 10, 1),
    @("SPVSX5", 1, 0.00000001, 1, 10, 1),
    @("PSC", 6.91, 0.0000283, 0.047600000000000003, -26000, 2.68)
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

$ws.Range("F9").Select()
